$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '66.365.96'
$ws.Range("E2").Value2 = '  +2.16%  '
$ws.Range("D3").Value2 = '3.417.44'
$ws.Range("E3").Value2 = '  +0.83%  '
$ws.Range("E4").Value2 = '  -0.06%  '
$ws.Range("D5").Value2 = '''567.69'
$ws.Range("E5").Value2 = '  +1.06%  '
$ws.Range("D6").Value2 = '''181.70'
$ws.Range("E6").Value2 = '  +4.26%  '
$ws.Range("E7").Value2 = '  +1.30%  '
$ws.Range("D8").Value2 = '3.411.28'
$ws.Range("E8").Value2 = '  +0.93%  '
$ws.Range("E9").Value2 = '  +0.05%  '
$ws.Range("D10").Value2 = '''0.178'
$ws.Range("E10").Value2 = '  +6.34%  '
$ws.Range("D11").Value2 = '''0.639'
$ws.Range("E11").Value2 = '  +1.30%  '
$ws.Range("D12").Value2 = '''54.77'
$ws.Range("E12").Value2 = '  +1.41%  '
$ws.Range("E13").Value2 = '  +0.47%  '
$ws.Range("E14").Value2 = '  +2.46%  '
$ws.Range("D15").Value2 = '3.965.34'
$ws.Range("E15").Value2 = '  +0.85%  '
$ws.Range("D16").Value2 = '''18.36'
$ws.Range("E16").Value2 = '  +0.36%  '
$ws.Range("B17").Value2 = 'WrappedEther'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value2 = '3.416.27'
$ws.Range("E17").Value2 = '  +0.85%  '
$ws.Range("B18").Value2 = 'TRON'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value2 = '''0.120'
$ws.Range("E18").Value2 = '  +0.50%  '
$ws.Range("D19").Value2 = '66.245.64'
$ws.Range("E19").Value2 = '  +1.98%  '
$ws.Range("D20").Value2 = '''11.99'
$ws.Range("E20").Value2 = '  +1.44%  '
$ws.Range("E21").Value2 = '  +1.24%  '
$ws.Range("D22").Value2 = '''471.46'
$ws.Range("E22").Value2 = '  +0.82%  '
$ws.Range("D23").Value2 = '''5.01'
$ws.Range("E23").Value2 = '  +2.37%  '
$ws.Range("E24").Value2 = '  +8.05%  '
$ws.Range("E25").Value2 = '  -0.03%  '
$ws.Range("D26").Value2 = '''89.45'
$ws.Range("E26").Value2 = '  +3.08%  '
$ws.Range("E27").Value2 = '  +0.47%  '
$ws.Range("D28").Value2 = '''10.82'
$ws.Range("E28").Value2 = '  +0.05%  '
$ws.Range("D29").Value2 = '''8.87'
$ws.Range("E29").Value2 = '  +1.19%  '
$ws.Range("D30").Value2 = '''31.37'
$ws.Range("E30").Value2 = '  +2.44%  '
$ws.Range("E31").Value2 = '  +4.01%  '
$ws.Range("E32").Value2 = '  +0.76%  '
$ws.Range("D33").Value2 = '''586.26'
$ws.Range("E33").Value2 = '  +2.39%  '
$ws.Range("D34").Value2 = '''62.62'
$ws.Range("E34").Value2 = '  +1.93%  '
$ws.Range("E35").Value2 = '  +1.19%  '
$ws.Range("E36").Value2 = '  -0.19%  '
$ws.Range("D37").Value2 = '''0.146'
$ws.Range("E37").Value2 = '  +4.61%  '
$ws.Range("D38").Value2 = '''3.60'
$ws.Range("E38").Value2 = '  -0.66%  '
$ws.Range("D39").Value2 = '''36.41'
$ws.Range("E39").Value2 = '  +2.34%  '
$ws.Range("D40").Value2 = '''0.385'
$ws.Range("E40").Value2 = '  +4.20%  '
$ws.Range("D42").Value2 = '3.131.73'
$ws.Range("E42").Value2 = '  +1.19%  '
$ws.Range("D43").Value2 = '''2.92'
$ws.Range("E43").Value2 = '  +2.56%  '
$ws.Range("E44").Value2 = '  +2.33%  '
$ws.Range("B45").Value2 = 'Fetch.AI'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value2 = '''2.53'
$ws.Range("E45").Value2 = '  +2.88%  '
$ws.Range("B46").Value2 = 'dogwifhat'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value2 = '''2.79'
$ws.Range("E46").Value2 = '  +19.50%  '
$ws.Range("E47").Value2 = '  +1.89%  '
$ws.Range("E48").Value2 = '  -0.12%  '
$ws.Range("E49").Value2 = '  -0.08%  '
$ws.Range("D50").Value2 = '''140.61'
$ws.Range("E50").Value2 = '  +1.87%  '
$ws.Range("D51").Value2 = '''8.59'
$ws.Range("E51").Value2 = '  +3.89%  '
